$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 - logistic_embeddings
$ws.Range("C5").Value = 0.274
$ws.Range("D5").Value = 0.403
$ws.Range("E5").Value = 0.443
$ws.Range("F5").Value = 0.488
$ws.Range("G5").Value = 0.5570000000000001
$ws.Range("H5").Value = 0.5610000000000001

# Row 7 - classical-best-embeddings -> classical-best-embed
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.274
$ws.Range("D7").Value = 0.403
$ws.Range("E7").Value = 0.443
$ws.Range("F7").Value = 0.488
$ws.Range("H7").Value = 0.5610000000000001

# Row 8 - BERT-base
$ws.Range("D8").Value = 0.508
$ws.Range("E8").Value = 0.549
$ws.Range("F8").Value = 0.576
$ws.Range("G8").Value = 0.64
$ws.Range("H8").Value = 0.65

# Row 9 - BERT-base-nli
$ws.Range("B9").Value = 0.256
$ws.Range("C9").Value = 0.389
$ws.Range("D9").Value = 0.527
$ws.Range("F9").Value = 0.589
$ws.Range("G9").Value = 0.624
$ws.Range("H9").Value = 0.631
